$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.983.96"
$ws.Range("E2").Value = "  +2.24%  "
$ws.Range("D3").Value = "2.050.92"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.33"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.79"
$ws.Range("E7").Value = "  +6.71%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +1.81%  "
$ws.Range("E10").Value = "  +2.85%  "
$ws.Range("E11").Value = "  +1.32%  "
$ws.Range("D12").Value = "2.354.09"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.63"
$ws.Range("E13").Value = "  +2.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.85"
$ws.Range("E14").Value = "  +3.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.28"
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("E16").Value = "  +1.03%  "
$ws.Range("D17").Value = "2.079.29"
$ws.Range("E17").Value = "  -4.71%  "
$ws.Range("D18").Value = "37.910.13"
$ws.Range("E18").Value = "  +2.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.27"
$ws.Range("E19").Value = "  -3.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.66"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("D21").Value = "0.0₃0835"
$ws.Range("E21").Value = "  +2.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.63"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.23"
$ws.Range("E25").Value = "  +2.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.28"
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.35"
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("E28").Value = "  +4.43%  "
$ws.Range("E29").Value = "  +1.62%  "
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("E31").Value = "  +1.39%  "
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.58"
$ws.Range("E33").Value = "  +2.53%  "
$ws.Range("E34").Value = "  +10.40%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.06"
$ws.Range("E37").Value = "  +9.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.28"
$ws.Range("E38").Value = "  +4.68%  "
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").Value = "1.487.63"
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("E41").Value = "  +1.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.00"
$ws.Range("E42").Value = "  +1.35%  "
$ws.Range("E43").Value = "  +2.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.53"
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0923"
$ws.Range("E45").Value = "  +1.56%  "
$ws.Range("E46").Value = "  -0.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.12"
$ws.Range("E47").Value = "  +12.52%  "
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.08"
$ws.Range("E50").Value = "  -2.34%  "
$ws.Range("D51").Value = "2.243.71"
$ws.Range("E51").Value = "  +1.50%  "
